# Agrega proceso de encapado y stamping
#
# On the "Maquinas" sheet, insert two new process rows:
#   - "Stamping"  (reuses existing shared strings "Stamping"/"Stamping")
#   - "Encapado"  (reuses existing shared strings "Encapado"/"Encapado")
# which pushes "Plastificado" down by one row and the rest of the table
# down by two rows overall. Also update the saved selection/active-sheet
# state: "OrdenEstandar" becomes the active/selected tab (cell F12),
# while "Maquinas" keeps cell A8 selected but is no longer the active tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Maquinas")

# Insert the "Stamping" row right above the current "Plastificado" row (row 6).
$ws.Range("A6").EntireRow.Insert()
$ws.Range("A6").Value = "Stamping"
$ws.Range("B6").Value = "Stamping"
$ws.Range("C6").Value = 1000
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = 10

# "Plastificado" is now on row 7. Insert the "Encapado" row right below it (row 8).
$ws.Range("A8").EntireRow.Insert()
$ws.Range("A8").Value = "Encapado"
$ws.Range("B8").Value = "Encapado"
$ws.Range("C8").Value = 1000
$ws.Range("D8").Value = 10
$ws.Range("E8").Value = 10

# Update the saved view state for the Maquinas sheet.
$ws.Range("A8").Select()

# OrdenEstandar becomes the active sheet with F12 selected.
$ws2 = $wb.Worksheets.Item("OrdenEstandar")
$ws2.Activate()
$ws2.Range("F12").Select()
